$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35 becomes the record previously in row 36
$ws.Range("A35").Value = 111966228
$ws.Range("B35").Value = 89183
$ws.Range("E35").Value = 3215
$ws.Range("F35").Value = "Rödgul trumpetsvamp"
$ws.Range("G35").Value = "Craterellus lutescens"
$ws.Range("H35").Value = "(Fr.) Fr."
$ws.Range("I35").Value = "'10"
$ws.Range("Q35").Value = 338356.4103134849
$ws.Range("R35").Value = 6433540.273063039
$ws.Range("AC35").Value = "I våtmarken öster om Angertuvan. Ca 35 m söder om hyggeskanten."
$ws.Range("AH35").Value = "Sumpskog"

# Row 36 becomes the record previously in row 35
$ws.Range("A36").Value = 111966065
$ws.Range("B36").Value = 83148
$ws.Range("E36").Value = 3518
$ws.Range("F36").Value = "Smal svampklubba"
$ws.Range("G36").Value = "Tolypocladium ophioglossoides"
$ws.Range("H36").Value = "(Ehrh. ex J.F.Gmel.:Fr.) Quandt, Kepler & Spatafora"
$ws.Range("I36").Value = "'2"
$ws.Range("Q36").Value = 338285.5070198396
$ws.Range("R36").Value = 6433442.904015562
$ws.Range("AC36").Value = "Växte vid största stigen som går vid foten av Angertuvans östra sluttning."
$ws.Range("AH36").ClearContents()
